$wb = $excel.ActiveWorkbook

# Update "想去人数" (F6 / F10) counts on the "展览" and "全部类型" sheets,
# which hold identical event listings.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F6").Value = 1612
    $ws.Range("F10").Value = 104
}
